$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell, well outside the used range, used to stage text-typed
# literals (numeric/date-looking strings, and the empty string) so that
# writing them into the real destination cells never lets Excel's
# type-inference turn them into numbers/dates. Set to Text ("@") format,
# write the literal, copy it, then PasteSpecial values-only into the real
# cell (which keeps the destination cell's own, unmodified, default style)
# and finally clear the scratch cell.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-TextValue($row, $col, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
}

# Row 2 - replace existing "Obra_008 / NAN-4224X" car entry with new data
$ws.Cells.Item(2, 1).Value = "CC_1609634842040"
$ws.Cells.Item(2, 2).Value = "Obra_067"
$ws.Cells.Item(2, 3).Value = "PHA-4543"
$ws.Cells.Item(2, 4).Value = "PÓ DE PEDRA"
Set-TextValue 2 5 "806"
$ws.Cells.Item(2, 6).Value = "null"
Set-TextValue 2 7 "-3.07019980852851"
Set-TextValue 2 8 "-60.0083330533303"
Set-TextValue 2 9 "2-1-2021"
Set-TextValue 2 10 "20:47"
Set-TextValue 2 11 ""

# Row 3 - replace existing "Obra_010 / JXL-1146" car entry with new data
$ws.Cells.Item(3, 1).Value = "CC_1609634868885"
$ws.Cells.Item(3, 2).Value = "Obra_067"
$ws.Cells.Item(3, 3).Value = "CC-11 A"
$ws.Cells.Item(3, 4).Value = "CONCRETO"
Set-TextValue 3 5 "465"
$ws.Cells.Item(3, 6).Value = "null"
Set-TextValue 3 7 "-3.07019980852851"
Set-TextValue 3 8 "-60.0083330533303"
Set-TextValue 3 9 "2-1-2021"
Set-TextValue 3 10 "20:47"
Set-TextValue 3 11 ""

# Row 4 - new row added
$ws.Cells.Item(4, 1).Value = "CC_1609634981898"
$ws.Cells.Item(4, 2).Value = "Obra_067"
$ws.Cells.Item(4, 3).Value = "NOZ-1549 N"
$ws.Cells.Item(4, 4).Value = "BRITA CONTAMINADA"
Set-TextValue 4 5 "943"
Set-TextValue 4 6 "287"
Set-TextValue 4 7 "-3.07021826877877"
Set-TextValue 4 8 "-60.0083191729802"
Set-TextValue 4 9 "2-1-2021"
Set-TextValue 4 10 "20:49"
Set-TextValue 4 11 "20:51"

# Remove the scratch cell and its staging value/format.
$scratch.Clear()
